$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Recorded By" values in column G for all recorded rows.
# Rows 125 and 141 are "Not Recorded" rows whose G cell is already blank,
# so they are naturally skipped by the three contiguous ranges below.
$ws.Range("G2:G124").ClearContents()
$ws.Range("G126:G140").ClearContents()
$ws.Range("G142:G153").ClearContents()

# Narrow column G from width 50 down to width 13.
# ColumnWidth is stored with a +5/6 character offset in the saved XML,
# so subtract that offset to land exactly on stored width 13.
$ws.Columns.Item(7).ColumnWidth = 13 - 5/6
